$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrB2 = New-Object 'object[,]' 24,3
$arrB2[0,0] = 0.7654999511490246
$arrB2[0,1] = 0.1755721224176696
$arrB2[0,2] = 0.06333682137619334
$arrB2[1,0] = 0.7271029672535008
$arrB2[1,1] = 0.1740247446825904
$arrB2[1,2] = 0.06252604219389468
$arrB2[2,0] = 0.7039131688004545
$arrB2[2,1] = 0.1730567633522355
$arrB2[2,2] = 0.06201887037882159
$arrB2[3,0] = 0.6945605978342826
$arrB2[3,1] = 0.1726578048892193
$arrB2[3,2] = 0.06180985158111341
$arrB2[4,0] = 0.6930135078772821
$arrB2[4,1] = 0.1725912864382977
$arrB2[4,2] = 0.06177500297357952
$arrB2[5,0] = 0.7037866416521297
$arrB2[5,1] = 0.1730514010706727
$arrB2[5,2] = 0.06201606094989032
$arrB2[6,0] = 0.7521807437069015
$arrB2[6,1] = 0.17504229775016
$arrB2[6,2] = 0.0630592105073795
$arrB2[7,0] = 0.850134373392649
$arrB2[7,1] = 0.1788047172473952
$arrB2[7,2] = 0.06503031878743215
$arrB2[8,0] = 0.9239561619995129
$arrB2[8,1] = 0.1814830441988704
$arrB2[8,2] = 0.06643276443971757
$arrB2[9,0] = 0.9579417340861482
$arrB2[9,1] = 0.1826829415226072
$arrB2[9,2] = 0.06706078661241577
$arrB2[10,0] = 0.9708689835372297
$arrB2[10,1] = 0.1831346570081038
$arrB2[10,2] = 0.0672971629368817
$arrB2[11,0] = 0.9680823113410497
$arrB2[11,1] = 0.1830374903006344
$arrB2[11,2] = 0.06724631936224768
$arrB2[12,0] = 0.959004113538839
$arrB2[12,1] = 0.1827201577363198
$arrB2[12,2] = 0.06708026239402898
$arrB2[13,0] = 0.9534509513431999
$arrB2[13,1] = 0.1825254359856672
$arrB2[13,2] = 0.06697835960834908
$arrB2[14,0] = 0.9217432139244011
$arrB2[14,1] = 0.1814042559797286
$arrB2[14,2] = 0.06639152047264218
$arrB2[15,0] = 0.9023946515559089
$arrB2[15,1] = 0.1807117125457509
$arrB2[15,2] = 0.06602895607582582
$arrB2[16,0] = 0.8913038892665099
$arrB2[16,1] = 0.180311639863362
$arrB2[16,2] = 0.06581948199119125
$arrB2[17,0] = 0.8875552896078034
$arrB2[17,1] = 0.1801758833130052
$arrB2[17,2] = 0.06574839710940239
$arrB2[18,0] = 0.9044504070955384
$arrB2[18,1] = 0.1807856150706471
$arrB2[18,2] = 0.06606764868743653
$arrB2[19,0] = 0.9616690388219524
$arrB2[19,1] = 0.1828134381920776
$arrB2[19,2] = 0.06712907656795863
$arrB2[20,0] = 0.999400571990634
$arrB2[20,1] = 0.1841232351641864
$arrB2[20,2] = 0.06781437104327637
$arrB2[21,0] = 0.979231943296611
$arrB2[21,1] = 0.183425590749863
$arrB2[21,2] = 0.06744938934599531
$arrB2[22,0] = 0.9035208972457269
$arrB2[22,1] = 0.1807522097208718
$arrB2[22,2] = 0.06605015896749222
$arrB2[23,0] = 0.8233090940973682
$arrB2[23,1] = 0.1778020063879566
$arrB2[23,2] = 0.06450509010410954
$ws.Range("B2:D25").Value = $arrB2

$arrF2 = New-Object 'object[,]' 24,3
$arrF2[0,0] = 1.353772499386906
$arrF2[0,1] = 1.240737194241873
$arrF2[0,2] = 1.18635307643612
$arrF2[1,0] = 1.348334223858807
$arrF2[1,1] = 1.236880471304218
$arrF2[1,2] = 1.189187681038732
$arrF2[2,0] = 1.345685181285155
$arrF2[2,1] = 1.235167028246693
$arrF2[2,2] = 1.191397437311863
$arrF2[3,0] = 1.344779225183188
$arrF2[3,1] = 1.234633330564122
$arrF2[3,2] = 1.19241602080119
$arrF2[4,0] = 1.344639274618132
$arrF2[4,1] = 1.23455464681598
$arrF2[4,2] = 1.192592290562729
$arrF2[5,0] = 1.345672260516665
$arrF2[5,1] = 1.23515916446172
$arrF2[5,2] = 1.191410696038233
$arrF2[6,0] = 1.351754164638749
$arrF2[6,1] = 1.239271492863907
$arrF2[6,2] = 1.187233092158209
$arrF2[7,0] = 1.3691576065942
$arrF2[7,1] = 1.252534531316655
$arrF2[7,2] = 1.182761759521028
$arrF2[8,0] = 1.385288169202695
$arrF2[8,1] = 1.265457884767756
$arrF2[8,2] = 1.181742625166237
$arrF2[9,0] = 1.393353960148943
$arrF2[9,1] = 1.272029635478361
$arrF2[9,2] = 1.18177061958329
$arrF2[10,0] = 1.396512985933839
$arrF2[10,1] = 1.274617947161033
$arrF2[10,2] = 1.181851861107489
$arrF2[11,0] = 1.39582797708087
$arrF2[11,1] = 1.274056070694513
$arrF2[11,2] = 1.18183122333636
$arrF2[12,0] = 1.393611756971026
$arrF2[12,1] = 1.272240578074204
$arrF2[12,2] = 1.181775887919301
$arrF2[13,0] = 1.392267891711413
$arrF2[13,1] = 1.271141527518466
$arrF2[13,2] = 1.181751191368591
$arrF2[14,0] = 1.384775699384519
$arrF2[14,1] = 1.265042355987347
$arrF2[14,2] = 1.181750683574307
$arrF2[15,0] = 1.380365927684849
$arrF2[15,1] = 1.261478239641335
$arrF2[15,2] = 1.181876246698351
$arrF2[16,0] = 1.377898057879577
$arrF2[16,1] = 1.259493462659819
$arrF2[16,2] = 1.181994742352032
$arrF2[17,0] = 1.377074246592215
$arrF2[17,1] = 1.258832647829379
$arrF2[17,2] = 1.18204281231823
$arrF2[18,0] = 1.380828264353255
$arrF2[18,1] = 1.261850896345834
$arrF2[18,2] = 1.181858091365811
$arrF2[19,0] = 1.394259873297486
$arrF2[19,1] = 1.272771124751344
$arrF2[19,2] = 1.181790224510365
$arrF2[20,0] = 1.403648401093463
$arrF2[20,1] = 1.280489488733878
$arrF2[20,2] = 1.182157594236941
$arrF2[21,0] = 1.398581730123396
$arrF2[21,1] = 1.276316824342942
$arrF2[21,2] = 1.181923865979826
$arrF2[22,0] = 1.380619032099688
$arrF2[22,1] = 1.261682218022599
$arrF2[22,2] = 1.181866155136618
$arrF2[23,0] = 1.36386262208238
$arrF2[23,1] = 1.248388964633094
$arrF2[23,2] = 1.183573276973647
$ws.Range("F2:H25").Value = $arrF2

$arrK2 = New-Object 'object[,]' 24,4
$arrK2[0,0] = 0.3856059181388503
$arrK2[0,1] = 0.2744899118084021
$arrK2[0,2] = 0.2262901638526174
$arrK2[0,3] = 2.395346053067556
$arrK2[1,0] = 0.3494102083500366
$arrK2[1,1] = 0.2711994070398376
$arrK2[1,2] = 0.2188805668913361
$arrK2[1,3] = 2.414386363511419
$arrK2[2,0] = 0.3273191877308363
$arrK2[2,1] = 0.269313512932932
$arrK2[2,2] = 0.2144494050240056
$arrK2[2,3] = 2.426767200345495
$arrK2[3,0] = 0.3183506570591135
$arrK2[3,1] = 0.2685788642012383
$arrK2[3,2] = 0.2126735114319835
$arrK2[3,3] = 2.431985976116053
$arrK2[4,0] = 0.3168634837530391
$arrK2[4,1] = 0.2684589233902486
$arrK2[4,2] = 0.2123804302781274
$arrK2[4,3] = 2.432863027214943
$arrK2[5,0] = 0.3271980980838407
$arrK2[5,1] = 0.2693034680040824
$arrK2[5,2] = 0.2144253337735371
$arrK2[5,3] = 2.426836880186563
$arrK2[6,0] = 0.3730980944102669
$arrK2[6,1] = 0.2733274608289022
$arrK2[6,2] = 0.2237108130191139
$arrK2[6,3] = 2.401767927409352
$arrK2[7,0] = 0.4641604835762507
$arrK2[7,1] = 0.2822841451164066
$arrK2[7,2] = 0.2428565180256115
$arrK2[7,3] = 2.358083257876672
$arrK2[8,0] = 0.5317061443498972
$arrK2[8,1] = 0.2895133864222998
$arrK2[8,2] = 0.2574928112322539
$arrK2[8,3] = 2.32932613291419
$arrK2[9,0] = 0.5625743694722871
$arrK2[9,1] = 0.2929429140789495
$arrK2[9,2] = 0.264274852356543
$arrK2[9,3] = 2.316968622617125
$arrK2[10,0] = 0.5742835701518914
$arrK2[10,1] = 0.2942618171124423
$arrK2[10,2] = 0.2668608024912231
$arrK2[10,3] = 2.312393343871754
$arrK2[11,0] = 0.5717608955484934
$arrK2[11,1] = 0.2939768694156868
$arrK2[11,2] = 0.2663030842311045
$arrK2[11,3] = 2.313374074249516
$arrK2[12,0] = 0.5635372916093218
$arrK2[12,1] = 0.2930510161499456
$arrK2[12,2] = 0.2644872449907112
$arrK2[12,3] = 2.316590121951783
$arrK2[13,0] = 0.5585027033862104
$arrK2[13,1] = 0.2924865355234942
$arrK2[13,2] = 0.2633772986812062
$arrK2[13,3] = 2.31857362234306
$arrK2[14,0] = 0.5296916561323997
$arrK2[14,1] = 0.2892920901419984
$arrK2[14,2] = 0.2570520755245269
$arrK2[14,3] = 2.330148307167207
$arrK2[15,0] = 0.512053037934038
$arrK2[15,1] = 0.2873684580211489
$arrK2[15,2] = 0.2532034367774116
$arrK2[15,3] = 2.33743456650835
$arrK2[16,0] = 0.5019210924621973
$arrK2[16,1] = 0.2862753016000141
$arrK2[16,2] = 0.2510014701402241
$arrK2[16,3] = 2.341693594058604
$arrK2[17,0] = 0.4984928854908617
$arrK2[17,1] = 0.2859074574902962
$arrK2[17,2] = 0.2502579286266169
$arrK2[17,3] = 2.343147331957013
$arrK2[18,0] = 0.5139293221025412
$arrK2[18,1] = 0.2875718593169267
$arrK2[18,2] = 0.2536119238746792
$arrK2[18,3] = 2.336651876364144
$arrK2[19,0] = 0.5659522206620977
$arrK2[19,1] = 0.2933224134350496
$arrK2[19,2] = 0.2650201202230988
$arrK2[19,3] = 2.315642660899769
$arrK2[20,0] = 0.6000690791790078
$arrK2[20,1] = 0.2971985415568525
$arrK2[20,2] = 0.2725793911845926
$arrK2[20,3] = 2.302519574186128
$arrK2[21,0] = 0.5818496586620938
$arrK2[21,1] = 0.2951190151416228
$arrK2[21,2] = 0.2685354357232796
$arrK2[21,3] = 2.309467980029659
$arrK2[22,0] = 0.5130810268169057
$arrK2[22,1] = 0.2874798618262133
$arrK2[22,2] = 0.2534272136257414
$arrK2[22,3] = 2.337005512341406
$arrK2[23,0] = 0.439412960259034
$arrK2[23,1] = 0.2797471313973858
$arrK2[23,2] = 0.2375769215921224
$arrK2[23,3] = 2.369314943023994
$ws.Range("K2:N25").Value = $arrK2
